$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every row of data
# (rows 2-375). The whole column was bumped from serial date 45171
# (2023-09-02) to 45172 (2023-09-03).
$ws.Range("C2:C375").Value = 45172
